$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 5, mirroring the formatting of row 4 ---
$ws.Range("A4:L4").Copy()
$ws.Range("A5:L5").PasteSpecial(-4122)
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight

# --- Cell values for the new "Submit" operator test row ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "google"
$ws.Range("C5").Value = "https://google.com"
$ws.Range("D5").Value = "SBMT"
$ws.Range("E5").Value = "//input[@id='lst-ib']"
$ws.Range("F5").Value = "What is the weather in Dallas today?"
$ws.Range("G5").Value = "SBMT"
$ws.Range("H5").Value = "//input[@id='lst-ib']"
$ws.Range("I5").Value = "What is the weather in Dallas today?"
$ws.Range("J5").Value = "FIND"
$ws.Range("K5").Value = "YES"
$ws.Range("L5").Value = "Dallas"

# --- Hyperlink on the url cell (matching rows 2-4) ---
$ws.Hyperlinks.Add($ws.Range("C5"), "https://google.com", "", "", "")

# restore the normal (non-hyperlink) cell formatting that Hyperlinks.Add overrides
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)

# apply the blue underlined "link" run formatting to the whole display text,
# split across two adjoining character ranges so the engine keeps it as rich text
# inside the shared string (matching the other url cells) instead of a cell-level style
$ws.Range("C5").Characters(1, 17).Font.Underline = $true
$ws.Range("C5").Characters(18, 1).Font.Underline = $true
$ws.Range("C5").Characters(1, 17).Font.Color = 16711680
$ws.Range("C5").Characters(18, 1).Font.Color = 16711680
$ws.Range("C5").Characters(1, 17).Font.Size = 10
$ws.Range("C5").Characters(18, 1).Font.Size = 10
$ws.Range("C5").Characters(1, 17).Font.Name = "Helvetica Neue"
$ws.Range("C5").Characters(18, 1).Font.Name = "Helvetica Neue"

# --- Column width tweaks for K and L (split from the shared 11-12 range) ---
$ws.Columns.Item(11).ColumnWidth = 9.285714285714286
$ws.Columns.Item(12).ColumnWidth = 13.142857142857142
